# Fixed update to excel issue
# Roll the 16-week forecast window forward by one week (the current
# "W1" becomes last week's data and drops off, a new "W16" week is
# appended) and refresh the forecast figures, then update the
# dependent Summary metrics.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Helper: write a value as plain text, avoiding Excel's automatic
# type-sniffing (which would otherwise turn date-looking or
# number-looking strings into real dates/numbers and leave a
# lingering number-format style on the cell).
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# New Week_Start_Date (col B) and forecast values (cols D-H) for rows 2..17
$rows = @(
    @{ Row=2;  Date="2025-02-02"; D=105; E=81; F=92;  G=102; H=116 },
    @{ Row=3;  Date="2025-02-09"; D=100; E=82; F=95;  G=105; H=121 },
    @{ Row=4;  Date="2025-02-16"; D=96;  E=86; F=100; G=111; H=128 },
    @{ Row=5;  Date="2025-02-23"; D=106; E=87; F=101; G=113; H=132 },
    @{ Row=6;  Date="2025-03-02"; D=113; E=87; F=102; G=116; H=137 },
    @{ Row=7;  Date="2025-03-09"; D=112; E=86; F=101; G=114; H=135 },
    @{ Row=8;  Date="2025-03-16"; D=109; E=84; F=99;  G=114; H=137 },
    @{ Row=9;  Date="2025-03-23"; D=95;  E=85; F=102; G=118; H=142 },
    @{ Row=10; Date="2025-03-30"; D=86;  E=84; F=100; G=115; H=137 },
    @{ Row=11; Date="2025-04-06"; D=94;  E=80; F=96;  G=112; H=137 },
    @{ Row=12; Date="2025-04-13"; D=105; E=81; F=97;  G=114; H=139 },
    @{ Row=13; Date="2025-04-20"; D=105; E=81; F=97;  G=114; H=139 },
    @{ Row=14; Date="2025-04-27"; D=103; E=80; F=95;  G=111; H=136 },
    @{ Row=15; Date="2025-05-04"; D=90;  E=75; F=90;  G=106; H=132 },
    @{ Row=16; Date="2025-05-11"; D=83;  E=75; F=90;  G=105; H=129 },
    @{ Row=17; Date="2025-05-18"; D=85;  E=74; F=89;  G=105; H=130 }
)

foreach ($r in $rows) {
    Set-TextValue $wsForecast.Cells.Item($r.Row, 2) $r.Date
    $wsForecast.Cells.Item($r.Row, 4).Value = $r.D
    $wsForecast.Cells.Item($r.Row, 5).Value = $r.E
    $wsForecast.Cells.Item($r.Row, 6).Value = $r.F
    $wsForecast.Cells.Item($r.Row, 7).Value = $r.G
    $wsForecast.Cells.Item($r.Row, 8).Value = $r.H
}

# Update the Summary sheet metrics that depend on the refreshed data
Set-TextValue $wsSummary.Cells.Item(2, 2)  "2022-12-25 to 2025-01-26"
Set-TextValue $wsSummary.Cells.Item(4, 2)  "304"
Set-TextValue $wsSummary.Cells.Item(6, 2)  "129"
Set-TextValue $wsSummary.Cells.Item(7, 2)  "66"
Set-TextValue $wsSummary.Cells.Item(8, 2)  "14326 units"
Set-TextValue $wsSummary.Cells.Item(9, 2)  "1588"
Set-TextValue $wsSummary.Cells.Item(10, 2) "836"
Set-TextValue $wsSummary.Cells.Item(11, 2) "407"
Set-TextValue $wsSummary.Cells.Item(12, 2) "113"
Set-TextValue $wsSummary.Cells.Item(14, 2) "83"
